$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Cells.Item(31, 1).Value = 112017392
$ws.Cells.Item(31, 2).Value = 90710
$ws.Cells.Item(31, 4).Value = "NT"
$ws.Cells.Item(31, 5).Value = 5449
$ws.Cells.Item(31, 6).Value = "Svart taggsvamp"
$ws.Cells.Item(31, 7).Value = "Phellodon niger"
$ws.Cells.Item(31, 8).Value = "(Fr.:Fr.) P.Karst."
$ws.Cells.Item(31, 9).Value = ""
$ws.Cells.Item(31, 10).Value = ""
$ws.Cells.Item(31, 17).Value = 682712.0453105029
$ws.Cells.Item(31, 18).Value = 6575457.539765021

# Row 32
$ws.Cells.Item(32, 1).Value = 112017326
$ws.Cells.Item(32, 2).Value = 90660
$ws.Cells.Item(32, 4).Value = "NT"
$ws.Cells.Item(32, 5).Value = 4362
$ws.Cells.Item(32, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(32, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(32, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(32, 9).Value = ""
$ws.Cells.Item(32, 10).Value = ""
$ws.Cells.Item(32, 17).Value = 682713.7813606198
$ws.Cells.Item(32, 18).Value = 6575496.010644327

# Row 33
$ws.Cells.Item(33, 1).Value = 112017447
$ws.Cells.Item(33, 2).Value = 90666
$ws.Cells.Item(33, 4).Value = "LC"
$ws.Cells.Item(33, 5).Value = 4364
$ws.Cells.Item(33, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(33, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(33, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(33, 9).Value = ""
$ws.Cells.Item(33, 10).Value = ""
$ws.Cells.Item(33, 17).Value = 682844.1942409466
$ws.Cells.Item(33, 18).Value = 6575513.554896963

# Row 34
$ws.Cells.Item(34, 1).Value = 112017224
$ws.Cells.Item(34, 2).Value = 90678
$ws.Cells.Item(34, 4).Value = "LC"
$ws.Cells.Item(34, 5).Value = 4366
$ws.Cells.Item(34, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(34, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(34, 8).Value = "Banker"
$ws.Cells.Item(34, 9).Value = ""
$ws.Cells.Item(34, 10).Value = ""
$ws.Cells.Item(34, 17).Value = 682702.748818734
$ws.Cells.Item(34, 18).Value = 6575490.872789856

# Row 35
$ws.Cells.Item(35, 1).Value = 112017465
$ws.Cells.Item(35, 2).Value = 88032
$ws.Cells.Item(35, 4).Value = "VU"
$ws.Cells.Item(35, 5).Value = 6276
$ws.Cells.Item(35, 6).Value = "Goliatmusseron"
$ws.Cells.Item(35, 7).Value = "Tricholoma matsutake"
$ws.Cells.Item(35, 8).Value = "(S.Ito & S.Imai) Singer"
$ws.Cells.Item(35, 9).NumberFormat = "@"
$ws.Cells.Item(35, 9).Value = "3"
$ws.Cells.Item(35, 10).Value = "fruktkroppar"
$ws.Cells.Item(35, 17).Value = 682896.4696766059
$ws.Cells.Item(35, 18).Value = 6575514.027787391

# Row 36
$ws.Cells.Item(36, 1).Value = 112017534
$ws.Cells.Item(36, 2).Value = 87992
$ws.Cells.Item(36, 4).Value = "VU"
$ws.Cells.Item(36, 5).Value = 1593
$ws.Cells.Item(36, 6).Value = "Lakritsmusseron"
$ws.Cells.Item(36, 7).Value = "Tricholoma apium"
$ws.Cells.Item(36, 8).Value = "Jul.Schäff."
$ws.Cells.Item(36, 9).NumberFormat = "@"
$ws.Cells.Item(36, 9).Value = "4"
$ws.Cells.Item(36, 10).Value = "fruktkroppar"
$ws.Cells.Item(36, 17).Value = 683072.5368938858
$ws.Cells.Item(36, 18).Value = 6575477.991881827

# Row 37
$ws.Cells.Item(37, 1).Value = 112017130
$ws.Cells.Item(37, 2).Value = 90666
$ws.Cells.Item(37, 4).Value = "LC"
$ws.Cells.Item(37, 5).Value = 4364
$ws.Cells.Item(37, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(37, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(37, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(37, 9).Value = ""
$ws.Cells.Item(37, 10).Value = ""
$ws.Cells.Item(37, 17).Value = 682695.3118543178
$ws.Cells.Item(37, 18).Value = 6575453.662799283

# Row 38
$ws.Cells.Item(38, 1).Value = 112017488
$ws.Cells.Item(38, 2).Value = 90678
$ws.Cells.Item(38, 4).Value = "LC"
$ws.Cells.Item(38, 5).Value = 4366
$ws.Cells.Item(38, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(38, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(38, 8).Value = "Banker"
$ws.Cells.Item(38, 9).Value = ""
$ws.Cells.Item(38, 10).Value = ""
$ws.Cells.Item(38, 17).Value = 682955.8308828628
$ws.Cells.Item(38, 18).Value = 6575473.896637772

# Row 39
$ws.Cells.Item(39, 1).Value = 112017159
$ws.Cells.Item(39, 2).Value = 90710
$ws.Cells.Item(39, 4).Value = "NT"
$ws.Cells.Item(39, 5).Value = 5449
$ws.Cells.Item(39, 6).Value = "Svart taggsvamp"
$ws.Cells.Item(39, 7).Value = "Phellodon niger"
$ws.Cells.Item(39, 8).Value = "(Fr.:Fr.) P.Karst."
$ws.Cells.Item(39, 9).Value = ""
$ws.Cells.Item(39, 10).Value = ""
$ws.Cells.Item(39, 17).Value = 682698.5384611045
$ws.Cells.Item(39, 18).Value = 6575482.480741166

# Row 40
$ws.Cells.Item(40, 1).Value = 112017512
$ws.Cells.Item(40, 2).Value = 88032
$ws.Cells.Item(40, 4).Value = "VU"
$ws.Cells.Item(40, 5).Value = 6276
$ws.Cells.Item(40, 6).Value = "Goliatmusseron"
$ws.Cells.Item(40, 7).Value = "Tricholoma matsutake"
$ws.Cells.Item(40, 8).Value = "(S.Ito & S.Imai) Singer"
$ws.Cells.Item(40, 9).NumberFormat = "@"
$ws.Cells.Item(40, 9).Value = "4"
$ws.Cells.Item(40, 10).Value = "fruktkroppar"
$ws.Cells.Item(40, 17).Value = 683036.8460961942
$ws.Cells.Item(40, 18).Value = 6575484.458868909

# Row 41
$ws.Cells.Item(41, 1).Value = 112017413
$ws.Cells.Item(41, 2).Value = 90709
$ws.Cells.Item(41, 4).Value = "NT"
$ws.Cells.Item(41, 5).Value = 5448
$ws.Cells.Item(41, 6).Value = "Svartvit taggsvamp"
$ws.Cells.Item(41, 7).Value = "Phellodon connatus"
$ws.Cells.Item(41, 8).Value = "(Schultz) nom.prov"
$ws.Cells.Item(41, 9).Value = ""
$ws.Cells.Item(41, 10).Value = ""
$ws.Cells.Item(41, 17).Value = 682733.9332997696
$ws.Cells.Item(41, 18).Value = 6575482.138353716

# Row 42
$ws.Cells.Item(42, 1).Value = 112017430
$ws.Cells.Item(42, 2).Value = 90709
$ws.Cells.Item(42, 4).Value = "NT"
$ws.Cells.Item(42, 5).Value = 5448
$ws.Cells.Item(42, 6).Value = "Svartvit taggsvamp"
$ws.Cells.Item(42, 7).Value = "Phellodon connatus"
$ws.Cells.Item(42, 8).Value = "(Schultz) nom.prov"
$ws.Cells.Item(42, 9).Value = ""
$ws.Cells.Item(42, 10).Value = ""
$ws.Cells.Item(42, 17).Value = 682793.1335561723
$ws.Cells.Item(42, 18).Value = 6575519.79500053

# Row 43
$ws.Cells.Item(43, 1).Value = 112017252
$ws.Cells.Item(43, 2).Value = 90666
$ws.Cells.Item(43, 4).Value = "LC"
$ws.Cells.Item(43, 5).Value = 4364
$ws.Cells.Item(43, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(43, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(43, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(43, 9).Value = ""
$ws.Cells.Item(43, 10).Value = ""
$ws.Cells.Item(43, 17).Value = 682710.810501094
$ws.Cells.Item(43, 18).Value = 6575493.820233095
